$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the credentials row (row 3)
$ws.Range("A3").Value = "Aman"
$ws.Range("B3").Value = "noPass"

# Remove the extra row (row 4: Hello / World)
$ws.Rows.Item(4).Delete()

# Update the active selection to B3
$ws.Range("B3").Select()
